# The deck ships two themes:
#   ppt/theme/theme1.xml -> linked from the (one) Slide Master -> "Integral" / "Red Violet"
#   ppt/theme/theme2.xml -> linked from the Notes Master        -> "Office Theme" / "Office"
#
# The authored change swaps the two themes' colour schemes (font/format
# schemes are already byte-identical between the two themes), i.e. the
# Slide Master ends up using the stock "Office" palette.
#
# Apply it through the real PowerPoint object model: walk the active
# presentation's Slide Master theme colour scheme and push each of the
# twelve MsoThemeColorSchemeIndex slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the corresponding "Office" theme RGB value.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function HexRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colour scheme (the twelve slots in MsoThemeColorSchemeIndex order).
$officeColors = @(
    (HexRGB 0x00 0x00 0x00),   # 1  dk1       000000
    (HexRGB 0xFF 0xFF 0xFF),   # 2  lt1       FFFFFF
    (HexRGB 0x44 0x54 0x6A),   # 3  dk2       44546A
    (HexRGB 0xE7 0xE6 0xE6),   # 4  lt2       E7E6E6
    (HexRGB 0x5B 0x9B 0xD5),   # 5  accent1   5B9BD5
    (HexRGB 0xED 0x7D 0x31),   # 6  accent2   ED7D31
    (HexRGB 0xA5 0xA5 0xA5),   # 7  accent3   A5A5A5
    (HexRGB 0xFF 0xC0 0x00),   # 8  accent4   FFC000
    (HexRGB 0x44 0x72 0xC4),   # 9  accent5   4472C4
    (HexRGB 0x70 0xAD 0x47),   # 10 accent6   70AD47
    (HexRGB 0x05 0x63 0xC1),   # 11 hlink     0563C1
    (HexRGB 0x95 0x4F 0x72)    # 12 folHlink  954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}

# Best-effort: also relabel the theme/colour-scheme to the stock "Office
# Theme" / "Office" names (harmless no-ops on hosts that keep these
# read-only, but correct for hosts that honour them). Deliberately NOT
# touching Design.Name / Master.Name here: on this host that writes to
# the Slide Master's <p:cSld name="..."> instead of the theme part,
# which is outside the scope of this change.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
